# daily auto push: 2026-01-18 18:40 UTC
# Insert two new log rows (2026/01/18 and 2026/01/19) immediately before the
# existing "2026/12/29" block, shifting the rest of the table down by two
# rows (old rows 665-706 become 667-708).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 665 down by inserting 2 blank rows at 665:666.
$ws.Range("A665:D666").EntireRow.Insert()

# Column A holds plain-text dates (e.g. "2026/12/29"), not real Excel date
# serials, in the source data. Force text formatting before assignment so
# the COM layer doesn't auto-coerce the "yyyy/mm/dd" strings into dates.
$ws.Range("A665:A666").NumberFormat = "@"

# New row 665: 2026/01/18 (Sunday), hour 22
$ws.Range("A665").Value = "2026/01/18"
$ws.Range("B665").Value = "日"
$ws.Range("C665").Value = 22
$ws.Range("D665").Value = 201

# New row 666: 2026/01/19 (Monday), hour 1
$ws.Range("A666").Value = "2026/01/19"
$ws.Range("B666").Value = "月"
$ws.Range("C666").Value = 1
$ws.Range("D666").Value = 201

# Dimension should now span down to row 708 (706 original + 2 inserted).
